# Update "想去人数" (interested-attendee counts) on the 展览 and 全部类型 sheets
# to reflect the latest scrape output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 110
$wsExpo.Range("F4").Value = 134
$wsExpo.Range("F5").Value = 2901
$wsExpo.Range("F6").Value = 291
$wsExpo.Range("F7").Value = 394

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 110
$wsAll.Range("F4").Value = 134
$wsAll.Range("F5").Value = 2901
$wsAll.Range("F6").Value = 291
$wsAll.Range("F9").Value = 394
